# Updates loading_percent values for rows 2-25 (columns C-I, K-O)
# per the "case with 380 kV done" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.52466335784678
$ws.Range("D2").Value = 4.891646116629243
$ws.Range("E2").Value = 12.91650814749353
$ws.Range("F2").Value = 23.77289540029187
$ws.Range("G2").Value = 27.89291468686337
$ws.Range("H2").Value = 14.08116195337141
$ws.Range("I2").Value = 22.59315978985907
$ws.Range("K2").Value = 12.00938450003745
$ws.Range("L2").Value = 9.621641979521394
$ws.Range("M2").Value = 15.50014742433758
$ws.Range("N2").Value = 17.41218633433042
$ws.Range("O2").Value = 21.31720966075225

$ws.Range("C3").Value = 10.49721249161125
$ws.Range("D3").Value = 4.829468571961828
$ws.Range("E3").Value = 12.93429230434142
$ws.Range("F3").Value = 23.79467360175473
$ws.Range("G3").Value = 27.92403406839404
$ws.Range("H3").Value = 14.12517015780638
$ws.Range("I3").Value = 22.64627156695885
$ws.Range("K3").Value = 11.58244689848209
$ws.Range("L3").Value = 9.646383559982615
$ws.Range("M3").Value = 15.3341072183209
$ws.Range("N3").Value = 17.44062827374131
$ws.Range("O3").Value = 21.38147067390466

$ws.Range("C4").Value = 10.482307526135
$ws.Range("D4").Value = 4.790333010745694
$ws.Range("E4").Value = 12.94747906818017
$ws.Range("F4").Value = 23.81461402007158
$ws.Range("G4").Value = 27.95265079953594
$ws.Range("H4").Value = 14.15452310439193
$ws.Range("I4").Value = 22.68392084557971
$ws.Range("K4").Value = 11.31273934784887
$ws.Range("L4").Value = 9.662731065436111
$ws.Range("M4").Value = 15.2331823697794
$ws.Range("N4").Value = 17.4597610181148
$ws.Range("O4").Value = 21.42572757055265

$ws.Range("C5").Value = 10.47672875231164
$ws.Range("D5").Value = 4.774152602589344
$ws.Range("E5").Value = 12.95342288231462
$ws.Range("F5").Value = 23.82438933386036
$ws.Range("G5").Value = 27.96669646489411
$ws.Range("H5").Value = 14.16707075910436
$ws.Range("I5").Value = 22.70052785841297
$ws.Range("K5").Value = 11.20108452868128
$ws.Range("L5").Value = 9.66968378806164
$ws.Range("M5").Value = 15.19235191944257
$ws.Range("N5").Value = 17.4679782474775
$ws.Range("O5").Value = 21.44496701602618

$ws.Range("C6").Value = 10.47583243611702
$ws.Range("D6").Value = 4.771452081116545
$ws.Range("E6").Value = 12.95444427813117
$ws.Range("F6").Value = 23.8261120550749
$ws.Range("G6").Value = 27.9691724937095
$ws.Range("H6").Value = 14.16918967882703
$ws.Range("I6").Value = 22.70336174762009
$ws.Range("K6").Value = 11.18244368706256
$ws.Range("L6").Value = 9.670855865080798
$ws.Range("M6").Value = 15.18559115545239
$ws.Range("N6").Value = 17.4693681294929
$ws.Range("O6").Value = 21.44823437254957

$ws.Range("C7").Value = 10.48223027812724
$ws.Range("D7").Value = 4.790115725330892
$ws.Range("E7").Value = 12.94755692031872
$ws.Range("F7").Value = 23.81473917855493
$ws.Range("G7").Value = 27.95283058164085
$ws.Range("H7").Value = 14.15468995372257
$ws.Range("I7").Value = 22.68413969640591
$ws.Range("K7").Value = 11.31124038736671
$ws.Range("L7").Value = 9.662823653742285
$ws.Range("M7").Value = 15.23263046053361
$ws.Range("N7").Value = 17.45987013507124
$ws.Range("O7").Value = 21.42598216754453

$ws.Range("C8").Value = 10.51479658416424
$ws.Range("D8").Value = 4.870411292444039
$ws.Range("E8").Value = 12.92216959117151
$ws.Range("F8").Value = 23.77904012166252
$ws.Range("G8").Value = 27.90166721575271
$ws.Range("H8").Value = 14.09585188969384
$ws.Range("I8").Value = 22.61042567574009
$ws.Range("K8").Value = 11.86383541840723
$ws.Range("L8").Value = 9.629933100085966
$ws.Range("M8").Value = 15.44270802562152
$ws.Range("N8").Value = 17.42164709214461
$ws.Range("O8").Value = 21.33836896507698

$ws.Range("C9").Value = 10.59390175427858
$ws.Range("D9").Value = 5.019901771545714
$ws.Range("E9").Value = 12.89037122430809
$ws.Range("F9").Value = 23.76121669521246
$ws.Range("G9").Value = 27.87703132899452
$ws.Range("H9").Value = 13.99898241660317
$ws.Range("I9").Value = 22.50594529229887
$ws.Range("K9").Value = 12.88137385761588
$ws.Range("L9").Value = 9.574595360273474
$ws.Range("M9").Value = 15.86104221844941
$ws.Range("N9").Value = 17.35990518478079
$ws.Range("O9").Value = 21.20476345246564

$ws.Range("C10").Value = 10.66097790777072
$ws.Range("D10").Value = 5.124405649927796
$ws.Range("E10").Value = 12.87796373176681
$ws.Range("F10").Value = 23.77996548976623
$ws.Range("G10").Value = 27.90530075579042
$ws.Range("H10").Value = 13.93911168773473
$ws.Range("I10").Value = 22.45372859235041
$ws.Range("K10").Value = 13.58157689103201
$ws.Range("L10").Value = 9.539506207542777
$ws.Range("M10").Value = 16.1699185653233
$ws.Range("N10").Value = 17.32255657418909
$ws.Range("O10").Value = 21.13004218513985

$ws.Range("C11").Value = 10.69335758755514
$ws.Range("D11").Value = 5.170696836608597
$ws.Range("E11").Value = 12.87469314243035
$ws.Range("F11").Value = 23.79539607897403
$ws.Range("G11").Value = 27.92824053869426
$ws.Range("H11").Value = 13.91433132793281
$ws.Range("I11").Value = 22.43532258306356
$ws.Range("K11").Value = 13.88860513967742
$ws.Range("L11").Value = 9.524748594381991
$ws.Range("M11").Value = 16.310268235385
$ws.Range("N11").Value = 17.30729683350597
$ws.Range("O11").Value = 21.10116739299875

$ws.Range("C12").Value = 10.70587991657755
$ws.Range("D12").Value = 5.188039121814883
$ws.Range("E12").Value = 12.87379530523187
$ws.Range("F12").Value = 23.80222883741178
$ws.Range("G12").Value = 27.93837448875182
$ws.Range("H12").Value = 13.90530094499
$ws.Range("I12").Value = 22.42912270661884
$ws.Range("K12").Value = 14.00312989060431
$ws.Range("L12").Value = 9.519333228659745
$ws.Range("M12").Value = 16.36335096919728
$ws.Range("N12").Value = 17.30176646288933
$ws.Range("O12").Value = 21.09097099309565

$ws.Range("C13").Value = 10.70317153167761
$ws.Range("D13").Value = 5.184312594520231
$ws.Range("E13").Value = 12.87397353156841
$ws.Range("F13").Value = 23.80071332858346
$ws.Range("G13").Value = 27.93612765603498
$ws.Range("H13").Value = 13.90723007639628
$ws.Range("I13").Value = 22.43042369748801
$ws.Range("K13").Value = 13.97854359664792
$ws.Range("L13").Value = 9.520491832539356
$ws.Range("M13").Value = 16.3519221828081
$ws.Range("N13").Value = 17.30294650056604
$ws.Range("O13").Value = 21.093134122009

$ws.Range("C14").Value = 10.6943826238279
$ws.Range("D14").Value = 5.172127391210435
$ws.Range("E14").Value = 12.87461245373931
$ws.Range("F14").Value = 23.79593840289729
$ws.Range("G14").Value = 27.92904528184711
$ws.Range("H14").Value = 13.91358130821502
$ws.Range("I14").Value = 22.43479707684194
$ws.Range("K14").Value = 13.89806249371398
$ws.Range("L14").Value = 9.524299602508824
$ws.Range("M14").Value = 16.31463691526716
$ws.Range("N14").Value = 17.30683687682188
$ws.Range("O14").Value = 21.10031372977041

$ws.Range("C15").Value = 10.68903290392849
$ws.Range("D15").Value = 5.164639002693479
$ws.Range("E15").Value = 12.87504815429512
$ws.Range("F15").Value = 23.79314238457747
$ws.Range("G15").Value = 27.92489548368283
$ws.Range("H15").Value = 13.91751765470897
$ws.Range("I15").Value = 22.43757621202204
$ws.Range("K15").Value = 13.84853648379106
$ws.Range("L15").Value = 9.526654498214285
$ws.Range("M15").Value = 16.29178897057372
$ws.Range("N15").Value = 17.30925214309209
$ws.Range("O15").Value = 21.10480759757263

$ws.Range("C16").Value = 10.65889874035339
$ws.Range("D16").Value = 5.121354630393463
$ws.Range("E16").Value = 12.87822518297599
$ws.Range("F16").Value = 23.77909571748006
$ws.Range("G16").Value = 27.90400426488951
$ws.Range("H16").Value = 13.94078057690104
$ws.Range("I16").Value = 22.45503917762787
$ws.Range("K16").Value = 13.56127283103806
$ws.Range("L16").Value = 9.540494875100393
$ws.Range("M16").Value = 16.16073959103676
$ws.Range("N16").Value = 17.32358859663608
$ws.Range("O16").Value = 21.13203236132787

$ws.Range("C17").Value = 10.64088519236241
$ws.Range("D17").Value = 5.094475807297783
$ws.Range("E17").Value = 12.8807817096379
$ws.Range("F17").Value = 23.77224447385952
$ws.Range("G17").Value = 27.89376882890994
$ws.Range("H17").Value = 13.9556807032418
$ws.Range("I17").Value = 22.46712268239208
$ws.Range("K17").Value = 13.38203567275333
$ws.Range("L17").Value = 9.549293886618671
$ws.Range("M17").Value = 16.08027425295058
$ws.Range("N17").Value = 17.33282625828255
$ws.Range("O17").Value = 21.15004586578145

$ws.Range("C18").Value = 10.63070042616475
$ws.Range("D18").Value = 5.078898777420467
$ws.Range("E18").Value = 12.88247562361221
$ws.Range("F18").Value = 23.76895367781025
$ws.Range("G18").Value = 27.88883070217314
$ws.Range("H18").Value = 13.96448193034802
$ws.Range("I18").Value = 22.47457609385934
$ws.Range("K18").Value = 13.27786542776314
$ws.Range("L18").Value = 9.554468241616799
$ws.Range("M18").Value = 16.03398013629895
$ws.Range("N18").Value = 17.33830242312105
$ws.Range("O18").Value = 21.16088821117295

$ws.Range("C19").Value = 10.62728251686472
$ws.Range("D19").Value = 5.073604800644433
$ws.Range("E19").Value = 12.88308755226658
$ws.Range("F19").Value = 23.76795115669028
$ws.Range("G19").Value = 27.88732177295031
$ws.Range("H19").Value = 13.96750154977813
$ws.Range("I19").Value = 22.47718609422509
$ws.Range("K19").Value = 13.24241289563009
$ws.Range("L19").Value = 9.556239671367033
$ws.Range("M19").Value = 16.01830486449327
$ws.Range("N19").Value = 17.34018455774319
$ws.Range("O19").Value = 21.16464186457059

$ws.Range("C20").Value = 10.64278458785403
$ws.Range("D20").Value = 5.097349274083229
$ws.Range("E20").Value = 12.88048643929577
$ws.Range("F20").Value = 23.77290656044795
$ws.Range("G20").Value = 27.89476019875345
$ws.Range("H20").Value = 13.95407064199788
$ws.Range("I20").Value = 22.46578427411203
$ws.Range("K20").Value = 13.40122798417929
$ws.Range("L20").Value = 9.548345482186853
$ws.Range("M20").Value = 16.08884153201238
$ws.Range("N20").Value = 17.33182603800437
$ws.Range("O20").Value = 21.14807845260511

$ws.Range("C21").Value = 10.69695712170326
$ws.Range("D21").Value = 5.175711620621206
$ws.Range("E21").Value = 12.87441554731386
$ws.Range("F21").Value = 23.79731408719365
$ws.Range("G21").Value = 27.93108629947428
$ws.Range("H21").Value = 13.91170620349581
$ws.Range("I21").Value = 22.43349160246418
$ws.Range("K21").Value = 13.92174960968952
$ws.Range("L21").Value = 9.523176473990022
$ws.Range("M21").Value = 16.32559058473059
$ws.Range("N21").Value = 17.30568744921288
$ws.Range("O21").Value = 21.09818486368362

$ws.Range("C22").Value = 10.73387866906767
$ws.Range("D22").Value = 5.225831561206126
$ws.Range("E22").Value = 12.87243310503557
$ws.Range("F22").Value = 23.81903146470341
$ws.Range("G22").Value = 27.96326019448707
$ws.Range("H22").Value = 13.88607864259338
$ws.Range("I22").Value = 22.41687536279984
$ws.Range("K22").Value = 14.25176479722157
$ws.Range("L22").Value = 9.507735485083698
$ws.Range("M22").Value = 16.47992425561837
$ws.Range("N22").Value = 17.29005054217405
$ws.Range("O22").Value = 21.06987815799397

$ws.Range("C23").Value = 10.71403675387379
$ws.Range("D23").Value = 5.199184200839468
$ws.Range("E23").Value = 12.87330978036536
$ws.Range("F23").Value = 23.80691416704091
$ws.Range("G23").Value = 27.94531804884606
$ws.Range("H23").Value = 13.89956794592212
$ws.Range("I23").Value = 22.42533275248876
$ws.Range("K23").Value = 14.07658643430019
$ws.Range("L23").Value = 9.515884427586416
$ws.Range("M23").Value = 16.39760322008507
$ws.Range("N23").Value = 17.29826414656828
$ws.Range("O23").Value = 21.08459172687928

$ws.Range("C24").Value = 10.6419253365083
$ws.Range("D24").Value = 5.096050565508049
$ws.Range("E24").Value = 12.88061923282465
$ws.Range("F24").Value = 23.77260521219031
$ws.Range("G24").Value = 27.89430905245192
$ws.Range("H24").Value = 13.95479781886609
$ws.Range("I24").Value = 22.46638779104682
$ws.Range("K24").Value = 13.39255463576993
$ws.Range("L24").Value = 9.548773895463951
$ws.Range("M24").Value = 16.08496836636975
$ws.Range("N24").Value = 17.33227772261017
$ws.Range("O24").Value = 21.14896640596335

$ws.Range("C25").Value = 10.57090271588754
$ws.Range("D25").Value = 4.980363224255687
$ws.Range("E25").Value = 12.89704848111566
$ws.Range("F25").Value = 23.7604447799968
$ws.Range("G25").Value = 27.87556220238534
$ws.Range("H25").Value = 14.02320487499586
$ws.Range("I25").Value = 22.52990687186422
$ws.Range("K25").Value = 12.61396546508221
$ws.Range("L25").Value = 9.574595360273474
$ws.Range("M25").Value = 15.86104221844941
$ws.Range("N25").Value = 17.37519777076676
$ws.Range("O25").Value = 21.2368012660667
